$wb = $excel.ActiveWorkbook

# 1) Update selection on PigeonHoleSort sheet (selecting it temporarily activates it)
$wsPigeon = $wb.Worksheets.Item("PigeonHoleSort")
$wsPigeon.Range("B33").Select()

# 2) Add the new QuickSort worksheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsQuick = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsQuick.Name = "QuickSort"

# 3) Populate cells in the same order the original author entered them
#    (this keeps the shared-strings table in the same append order)
$wsQuick.Range("C1").Value = 'Temporal'
$wsQuick.Range("B2").Value = 'public void quickSort(int[] list, int high, int low ) {'
$wsQuick.Range("F2").Value = 'private int partition(int[] list, int low, int high) {'
$wsQuick.Range("B4").Value = '             if(low<high)'
$wsQuick.Range("F4").Value = 'int pivot = list[high];'
$wsQuick.Range("B5").Value = '              {'
$wsQuick.Range("F5").Value = 'int smallerElementIndex = (low-1);'
$wsQuick.Range("B6").Value = '                      int partition = partition(list,low,high);'
$wsQuick.Range("C6").Value = 'n-2'
$wsQuick.Range("F7").Value = 'for (int i = low; i < high; i++) {'
$wsQuick.Range("B8").Value = '                     quickSort(list,low, partition-1);'
$wsQuick.Range("B9").Value = '                     quickSort(list,partition+1,high);'
$wsQuick.Range("F9").Value = 'if(list[i]<=pivot)'
$wsQuick.Range("B10").Value = '               }'
$wsQuick.Range("F10").Value = '{'
$wsQuick.Range("F11").Value = '            smallerElementIndex++;'
$wsQuick.Range("F14").Value = '            int temp = list[smallerElementIndex];'
$wsQuick.Range("B15").Value = 'T(n) = [(n-1)n]/2'
$wsQuick.Range("F15").Value = '            list[smallerElementIndex] = list[i];'
$wsQuick.Range("F16").Value = '            list[i] = temp;'
$wsQuick.Range("F18").Value = '     }'
$wsQuick.Range("F22").Value = 'int temp = list[smallerElementIndex+1];'
$wsQuick.Range("F23").Value = 'list[smallerElementIndex+1] = list[high];'
$wsQuick.Range("F24").Value = 'list[high] = temp;'
$wsQuick.Range("F26").Value = 'return smallerElementIndex+1;'
$wsQuick.Range("D1").Value = 'Space'
$wsQuick.Range("G1").Value = 'Temporal'
$wsQuick.Range("H1").Value = 'Space'
$wsQuick.Range("C4").Value = 'n-1'
$wsQuick.Range("G7").Value = 'n+1'
$wsQuick.Range("G9").Value = 'n'
$wsQuick.Range("G11").Value = 'n'
$wsQuick.Range("B12").Value = '}'
$wsQuick.Range("G14").Value = 'n'
$wsQuick.Range("G15").Value = 'n'
$wsQuick.Range("G16").Value = 'n'
$wsQuick.Range("F19").Value = '}'
$wsQuick.Range("F27").Value = '}'
$wsQuick.Range("G4").Value = 1
$wsQuick.Range("G5").Value = 1
$wsQuick.Range("G22").Value = 1
$wsQuick.Range("G23").Value = 1
$wsQuick.Range("G24").Value = 1
$wsQuick.Range("G26").Value = 1

# 4) Set column widths to match the source sheet
$wsQuick.Columns.Item(2).ColumnWidth = 45.140625
$wsQuick.Columns.Item(6).ColumnWidth = 43.5703125

# 5) Final selection + activate QuickSort as the active tab
$wsQuick.Range("B15").Select()
